# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Set-NumericLookingText($cellRef, $val) {
    # Force Excel to store a numeric-looking string as text, matching the
    # original workbook (all D/E columns are plain text cells), then reset
    # the cell style back to Normal so no stray formatting is introduced.
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-PlainText "D2" "63.726.86"
Set-PlainText "E2" "  +2.75%  "
Set-PlainText "D3" "2.560.54"
Set-PlainText "E3" "  +5.87%  "
Set-PlainText "E4" "  +0.06%  "
Set-NumericLookingText "D5" "574.51"
Set-PlainText "E5" "  +2.96%  "
Set-NumericLookingText "D6" "151.19"
Set-PlainText "E6" "  +9.01%  "
Set-NumericLookingText "D8" "0.588"
Set-PlainText "E8" "  +0.91%  "
Set-PlainText "D9" "2.559.06"
Set-PlainText "E9" "  +6.04%  "
Set-PlainText "E10" "  +2.62%  "
Set-PlainText "E11" "  +0.09%  "
Set-PlainText "E12" "  +1.64%  "
Set-NumericLookingText "D13" "0.360"
Set-PlainText "E13" "  +4.14%  "
Set-NumericLookingText "D14" "28.27"
Set-PlainText "E14" "  +9.78%  "
Set-PlainText "D15" "3.019.87"
Set-PlainText "E15" "  +6.04%  "
Set-PlainText "D16" "63.642.93"
Set-PlainText "E16" "  +2.72%  "
Set-PlainText "E17" "  +3.65%  "
Set-PlainText "D18" "2.561.78"
Set-PlainText "E18" "  +6.17%  "
Set-NumericLookingText "D19" "11.66"
Set-PlainText "E19" "  +5.32%  "
Set-NumericLookingText "D20" "343.09"
Set-PlainText "E20" "  +0.01%  "
Set-PlainText "E21" "  +3.97%  "
Set-PlainText "E22" "  +1.16%  "
Set-PlainText "E23" "  -0.07%  "
Set-NumericLookingText "D24" "66.29"
Set-PlainText "E24" "  +2.12%  "
Set-PlainText "E25" "  -0.15%  "
Set-NumericLookingText "D26" "1.58"
Set-PlainText "E26" "  +3.85%  "
Set-PlainText "E27" "  -0.01%  "
Set-NumericLookingText "D28" "8.46"
Set-PlainText "E28" "  +2.25%  "
Set-NumericLookingText "D29" "1.44"
Set-PlainText "E29" "  +5.69%  "
Set-PlainText "E30" "  +13.09%  "
Set-PlainText "D31" "0.0₃0844"
Set-PlainText "E31" "  +7.61%  "
Set-PlainText "E32" "  +4.72%  "
Set-NumericLookingText "D33" "177.45"
Set-PlainText "E33" "  +3.71%  "
Set-NumericLookingText "D34" "1.59"
Set-PlainText "E34" "  +10.82%  "
Set-NumericLookingText "D35" "420.54"
Set-PlainText "E35" "  +11.95%  "
Set-PlainText "E36" "  +3.32%  "
Set-NumericLookingText "D37" "19.19"
Set-PlainText "E37" "  +3.76%  "
Set-NumericLookingText "D38" "4.48"
Set-PlainText "E38" "  +0.12%  "
Set-PlainText "E39" "  -0.02%  "
Set-PlainText "E40" "  +5.78%  "
Set-NumericLookingText "D41" "1.00"
Set-PlainText "E41" "  +0.16%  "
Set-NumericLookingText "D42" "40.51"
Set-PlainText "E42" "  +3.73%  "
Set-NumericLookingText "D43" "156.80"
Set-PlainText "E43" "  +7.66%  "
Set-NumericLookingText "D44" "3.82"
Set-PlainText "E44" "  +4.62%  "
Set-PlainText "E45" "  +2.39%  "
Set-PlainText "E46" "  +4.82%  "
Set-NumericLookingText "D47" "0.0534"
Set-PlainText "E47" "  +3.66%  "
Set-NumericLookingText "D48" "0.0970"
Set-PlainText "E48" "  +1.52%  "
# Row 49 <-> Row 50 swap: VeChain moves to rank 49 (was EnergySwap),
# EnergySwap moves to rank 50 (was VeChain), with updated price/volume figures.
Set-PlainText "B49" "VeChain"
Set-PlainText "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-NumericLookingText "D49" "0.0234"
Set-PlainText "E49" "  +6.08%  "

Set-PlainText "B50" "EnergySwap"
Set-PlainText "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-NumericLookingText "D50" "18.90"
Set-PlainText "E50" "  +4.88%  "

Set-NumericLookingText "D51" "1.88"
Set-PlainText "E51" "  +10.85%  "

